$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 118.0346986666667
$ws.Range("H2").Value = 354.104096
$ws.Range("I2").Value = 0.2666057129183408
$ws.Range("J2").Value = 0.2666057129183408
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.103724333333333
$ws.Range("N2").Value = 3.311173
$ws.Range("O2").Value = 0.01870879385910814
$ws.Range("P2").Value = 0.01870879385910814
$ws.Range("Q2").Value = 130.2777690960675
$ws.Range("R2").Value = 1172.499921864608
$ws.Range("S2").Value = 0.004987871324649804
$ws.Range("T2").Value = 0.004987871324649804

$ws.Range("G3").Value = 118.0346986666667
$ws.Range("H3").Value = 354.104096
$ws.Range("I3").Value = 0.2666057129183408
$ws.Range("J3").Value = 0.2666057129183408
$ws.Range("O3").Value = 0.1603368629650925
$ws.Range("P3").Value = 0.1603368629650925
$ws.Range("Q3").Value = 1116.497886943413
$ws.Range("R3").Value = 10048.48098249072
$ws.Range("S3").Value = 0.04274672365789881
$ws.Range("T3").Value = 0.04274672365789881

$ws.Range("G4").Value = 118.0346986666667
$ws.Range("H4").Value = 354.104096
$ws.Range("I4").Value = 0.2666057129183408
$ws.Range("J4").Value = 0.2666057129183408
$ws.Range("M4").Value = 47.61312599999999
$ws.Range("N4").Value = 142.839378
$ws.Range("O4").Value = 0.80707123365805
$ws.Range("P4").Value = 0.80707123365805
$ws.Range("Q4").Value = 5620.000979988031
$ws.Range("R4").Value = 50580.00881989227
$ws.Range("S4").Value = 0.2151698016252893
$ws.Range("T4").Value = 0.2151698016252893

$ws.Range("G5").Value = 118.0346986666667
$ws.Range("H5").Value = 354.104096
$ws.Range("I5").Value = 0.2666057129183408
$ws.Range("J5").Value = 0.2666057129183408
$ws.Range("M5").Value = 0.8190333333333334
$ws.Range("N5").Value = 2.4571
$ws.Range("O5").Value = 0.01388310951774934
$ws.Range("P5").Value = 0.01388310951774934
$ws.Range("Q5").Value = 96.67435269795556
$ws.Range("R5").Value = 870.0691742816
$ws.Range("S5").Value = 0.003701316310502965
$ws.Range("T5").Value = 0.003701316310502965

$ws.Range("I6").Value = 0.4881754016778185
$ws.Range("J6").Value = 0.4881754016778186
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.103724333333333
$ws.Range("N6").Value = 3.311173
$ws.Range("O6").Value = 0.01870879385910814
$ws.Range("P6").Value = 0.01870879385910814
$ws.Range("Q6").Value = 238.548535070749
$ws.Range("R6").Value = 2146.936815636741
$ws.Range("S6").Value = 0.009133172957077622
$ws.Range("T6").Value = 0.009133172957077624

$ws.Range("I7").Value = 0.4881754016778185
$ws.Range("J7").Value = 0.4881754016778186
$ws.Range("O7").Value = 0.1603368629650925
$ws.Range("P7").Value = 0.1603368629650925
$ws.Range("S7").Value = 0.0782725124817454
$ws.Range("T7").Value = 0.0782725124817454

$ws.Range("I8").Value = 0.4881754016778185
$ws.Range("J8").Value = 0.4881754016778186
$ws.Range("M8").Value = 47.61312599999999
$ws.Range("N8").Value = 142.839378
$ws.Range("O8").Value = 0.80707123365805
$ws.Range("P8").Value = 0.80707123365805
$ws.Range("Q8").Value = 10290.65058585491
$ws.Range("R8").Value = 92615.85527269421
$ws.Range("S8").Value = 0.3939923236736311
$ws.Range("T8").Value = 0.3939923236736311

$ws.Range("I9").Value = 0.4881754016778185
$ws.Range("J9").Value = 0.4881754016778186
$ws.Range("M9").Value = 0.8190333333333334
$ws.Range("N9").Value = 2.4571
$ws.Range("O9").Value = 0.01388310951774934
$ws.Range("P9").Value = 0.01388310951774934
$ws.Range("Q9").Value = 177.0181157923
$ws.Range("R9").Value = 1593.1630421307
$ws.Range("S9").Value = 0.006777392565364427
$ws.Range("T9").Value = 0.006777392565364428

$ws.Range("G10").Value = 45.876452
$ws.Range("H10").Value = 137.629356
$ws.Range("I10").Value = 0.1036214293744632
$ws.Range("J10").Value = 0.1036214293744632
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.103724333333333
$ws.Range("N10").Value = 3.311173
$ws.Range("O10").Value = 0.01870879385910814
$ws.Range("P10").Value = 0.01870879385910814
$ws.Range("Q10").Value = 50.63495639939867
$ws.Range("R10").Value = 455.714607594588
$ws.Range("S10").Value = 0.001938631961552965
$ws.Range("T10").Value = 0.001938631961552965

$ws.Range("G11").Value = 45.876452
$ws.Range("H11").Value = 137.629356
$ws.Range("I11").Value = 0.1036214293744632
$ws.Range("J11").Value = 0.1036214293744632
$ws.Range("O11").Value = 0.1603368629650925
$ws.Range("P11").Value = 0.1603368629650925
$ws.Range("Q11").Value = 433.9483414373801
$ws.Range("R11").Value = 3905.53507293642
$ws.Range("S11").Value = 0.01661433492186032
$ws.Range("T11").Value = 0.01661433492186032

$ws.Range("G12").Value = 45.876452
$ws.Range("H12").Value = 137.629356
$ws.Range("I12").Value = 0.1036214293744632
$ws.Range("J12").Value = 0.1036214293744632
$ws.Range("M12").Value = 47.61312599999999
$ws.Range("N12").Value = 142.839378
$ws.Range("O12").Value = 0.80707123365805
$ws.Range("P12").Value = 0.80707123365805
$ws.Range("Q12").Value = 2184.321289508952
$ws.Range("R12").Value = 19658.89160558057
$ws.Range("S12").Value = 0.0836298748386585
$ws.Range("T12").Value = 0.08362987483865851

$ws.Range("G13").Value = 45.876452
$ws.Range("H13").Value = 137.629356
$ws.Range("I13").Value = 0.1036214293744632
$ws.Range("J13").Value = 0.1036214293744632
$ws.Range("M13").Value = 0.8190333333333334
$ws.Range("N13").Value = 2.4571
$ws.Range("O13").Value = 0.01388310951774934
$ws.Range("P13").Value = 0.01388310951774934
$ws.Range("Q13").Value = 37.57434340306667
$ws.Range("R13").Value = 338.1690906276
$ws.Range("S13").Value = 0.001438587652391401
$ws.Range("T13").Value = 0.001438587652391401

$ws.Range("G14").Value = 62.68962833333333
$ws.Range("H14").Value = 188.068885
$ws.Range("I14").Value = 0.1415974560293775
$ws.Range("J14").Value = 0.1415974560293775
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.103724333333333
$ws.Range("N14").Value = 3.311173
$ws.Range("O14").Value = 0.01870879385910814
$ws.Range("P14").Value = 0.01870879385910814
$ws.Range("Q14").Value = 69.19206823912278
$ws.Range("R14").Value = 622.728614152105
$ws.Range("S14").Value = 0.002649117615827752
$ws.Range("T14").Value = 0.002649117615827753

$ws.Range("G15").Value = 62.68962833333333
$ws.Range("H15").Value = 188.068885
$ws.Range("I15").Value = 0.1415974560293775
$ws.Range("J15").Value = 0.1415974560293775
$ws.Range("O15").Value = 0.1603368629650925
$ws.Range("P15").Value = 0.1603368629650925
$ws.Range("Q15").Value = 592.9852692308417
$ws.Range("R15").Value = 5336.867423077575
$ws.Range("S15").Value = 0.02270329190358801
$ws.Range("T15").Value = 0.02270329190358801

$ws.Range("G16").Value = 62.68962833333333
$ws.Range("H16").Value = 188.068885
$ws.Range("I16").Value = 0.1415974560293775
$ws.Range("J16").Value = 0.1415974560293775
$ws.Range("M16").Value = 47.61312599999999
$ws.Range("N16").Value = 142.839378
$ws.Range("O16").Value = 0.80707123365805
$ws.Range("P16").Value = 0.80707123365805
$ws.Range("Q16").Value = 2984.849172728169
$ws.Range("R16").Value = 26863.64255455353
$ws.Range("S16").Value = 0.1142792335204712
$ws.Range("T16").Value = 0.1142792335204712

$ws.Range("G17").Value = 62.68962833333333
$ws.Range("H17").Value = 188.068885
$ws.Range("I17").Value = 0.1415974560293775
$ws.Range("J17").Value = 0.1415974560293775
$ws.Range("M17").Value = 0.8190333333333334
$ws.Range("N17").Value = 2.4571
$ws.Range("O17").Value = 0.01388310951774934
$ws.Range("P17").Value = 0.01388310951774934
$ws.Range("Q17").Value = 51.34489525927778
$ws.Range("R17").Value = 462.1040573335
$ws.Range("S17").Value = 0.001965812989490543
$ws.Range("T17").Value = 0.001965812989490544
